$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.45
$ws.Range("I2").Value = 6.25
$ws.Range("J2").Value = 1.92
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 12
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 11
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 201
$ws.Range("AJ2").Value = 19
$ws.Range("AU2").Value = 8
$ws.Range("AX2").Value = 8
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 101
$ws.Range("H4").Value = 3.3
$ws.Range("J4").Value = 2.75
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3.2
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.73
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.91
$ws.Range("AA4").Value = 17
$ws.Range("AC4").Value = 8.5
$ws.Range("AG4").Value = 301
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 19
$ws.Range("AL4").Value = 34
$ws.Range("AT4").Value = 2.63
$ws.Range("AV4").Value = 51
$ws.Range("O5").Value = 1.02
$ws.Range("P5").Value = 15
$ws.Range("Q5").Value = 1.05
$ws.Range("R5").Value = 9
$ws.Range("U5").Value = 1.91
$ws.Range("V5").Value = 1.85
$ws.Range("AB5").Value = 300
$ws.Range("AC5").Value = 50
$ws.Range("AF5").Value = 120
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 24
$ws.Range("AI5").Value = 11.75
$ws.Range("AK5").Value = 8.25
$ws.Range("AO5").Value = 300
$ws.Range("AP5").Value = 110
$ws.Range("AT5").Value = 7.5
$ws.Range("AU5").Value = 12.5
$ws.Range("AV5").Value = 65
$ws.Range("AX5").Value = 4.1
$ws.Range("AY5").Value = 4.05
$ws.Range("AZ5").Value = 12
$ws.Range("BA5").Value = 5.8
$ws.Range("BB5").Value = 17
$ws.Range("BC5").Value = 90
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("O7").Value = 1.36
$ws.Range("P7").Value = 3
$ws.Range("V7").Value = 1.54
$ws.Range("M8").Value = 1.11
$ws.Range("O8").Value = 1.5
$ws.Range("V8").Value = 1.54
$ws.Range("G9").Value = 2.6
$ws.Range("I9").Value = 3.1
$ws.Range("M9").Value = 1.17
$ws.Range("N9").Value = 5
$ws.Range("O9").Value = 1.67
$ws.Range("Q9").Value = 3.4
$ws.Range("R9").Value = 1.33
$ws.Range("S9").Value = 1.73
$ws.Range("T9").Value = 2.08
$ws.Range("U9").Value = 2.5
$ws.Range("V9").Value = 1.47
$ws.Range("W9").Value = 5.5
$ws.Range("AE9").Value = 23
$ws.Range("AH9").Value = 6
$ws.Range("AP9").Value = 41
$ws.Range("AT9").Value = 2
$ws.Range("AV9").Value = 101
$ws.Range("BA9").Value = 81
$ws.Range("G15").Value = 1.42
$ws.Range("H15").Value = 4.1
$ws.Range("I15").Value = 7.5
$ws.Range("J15").Value = 2
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("W15").Value = 6.5
$ws.Range("AC15").Value = 10
$ws.Range("AH15").Value = 19
$ws.Range("AI15").Value = 41
$ws.Range("AO15").Value = 7
$ws.Range("AV15").Value = 67
$ws.Range("BA15").Value = 151
$ws.Range("I17").Value = 3.25
$ws.Range("AB17").Value = 29
$ws.Range("AK17").Value = 34
$ws.Range("AO17").Value = 12
$ws.Range("G21").Value = 2.1
$ws.Range("J21").Value = 2.87
$ws.Range("M30").Value = 1.1
$ws.Range("N30").Value = 7
$ws.Range("O31").Value = 1.13
$ws.Range("P31").Value = 5.98
$ws.Range("U31").Value = 1.38
$ws.Range("V31").Value = 2.71
$ws.Range("R61").Value = 1.54
$ws.Range("K62").Value = 1.95
$ws.Range("R62").Value = 1.5
$ws.Range("G63").Value = 2.35
$ws.Range("I63").Value = 2.75
$ws.Range("J63").Value = 3.2
$ws.Range("L63").Value = 3.6
$ws.Range("M63").Value = 1.07
$ws.Range("N63").Value = 9
$ws.Range("R63").Value = 1.67
$ws.Range("Y63").Value = 10
$ws.Range("Z63").Value = 23
$ws.Range("AE63").Value = 15
$ws.Range("AH63").Value = 8
$ws.Range("AI63").Value = 13
$ws.Range("AK63").Value = 29
$ws.Range("AL63").Value = 23
$ws.Range("AN63").Value = 4.5
$ws.Range("AO63").Value = 15
$ws.Range("AP63").Value = 26
$ws.Range("AQ63").Value = 51
$ws.Range("AX63").Value = 4.75
$ws.Range("AZ63").Value = 26
